$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cells = [ordered]@{
    'D2' = '24.855.20'
    'E2' = '  +1.76%  '
    'D3' = '1.709.59'
    'E3' = '  +1.75%  '
    'E4' = '  +0.20%  '
    'D5' = '311.17'
    'D6' = '0.9993'
    'D7' = '0.3755'
    'E7' = '  +1.17%  '
    'D8' = '49.69'
    'E8' = '  +3.13%  '
    'D9' = '0.3446'
    'E9' = '  +0.12%  '
    'D10' = '1.208'
    'E10' = '  +2.03%  '
    'D11' = '0.07548'
    'E11' = '  +3.72%  '
    'D12' = '1.002'
    'E12' = '  +0.25%  '
    'D13' = '21.11'
    'E13' = '  +3.18%  '
    'D14' = '6.298'
    'E14' = '  +2.77%  '
    'D15' = '7.042'
    'D16' = '1.707.87'
    'E16' = '  +1.81%  '
    'D17' = '0.00001137'
    'E17' = '  +2.26%  '
    'D18' = '0.06740'
    'E18' = '  +0.18%  '
    'D19' = '0.9994'
    'E19' = '  +0.33%  '
    'D20' = '84.61'
    'E20' = '  +4.09%  '
    'D21' = '17.33'
    'E21' = '  +5.35%  '
    'D22' = '6.383'
    'E22' = '  +4.56%  '
    'D23' = '13.26'
    'E23' = '  +10.79%  '
    'D24' = '24.819.78'
    'E24' = '  +1.80%  '
    'D25' = '2.447'
    'E25' = '  +0.65%  '
    'D26' = '2.794'
    'E26' = '  +4.20%  '
    'D27' = '20.39'
    'E27' = '  +3.91%  '
    'D28' = '152.00'
    'E28' = '  -0.34%  '
    'D29' = '132.26'
    'E29' = '  +3.91%  '
    'D30' = '1.898.45'
    'E30' = '  +1.98%  '
    'D31' = '1.239'
    'E31' = '  +27.22%  '
    'D32' = '6.951'
    'E32' = '  +9.72%  '
    'D33' = '4.256'
    'E33' = '  +5.91%  '
    'D34' = '1.832'
    'E34' = '  +5.77%  '
    'D35' = '13.78'
    'E35' = '  +11.62%  '
    'D36' = '0.08783'
    'E36' = '  +3.47%  '
    'D37' = '5.617'
    'E37' = '  +4.99%  '
    'D38' = '9.336'
    'E38' = '  +2.82%  '
    'D39' = '0.06693'
    'E39' = '  +2.97%  '
    'D40' = '0.02410'
    'E40' = '  +3.10%  '
    'D41' = '0.2239'
    'E41' = '  +5.76%  '
    'D42' = '1.280'
    'E42' = '  +1.35%  '
    'D43' = '0.6451'
    'E43' = '  +4.11%  '
    'D44' = '0.9995'
    'E44' = '  +0.33%  '
    'D45' = '13.98'
    'E45' = '  +7.63%  '
    'D46' = '0.6170'
    'E46' = '  +3.54%  '
    'D47' = '3.826'
    'E47' = '  +1.11%  '
    'E48' = '  +5.36%  '
    'D49' = '130.36'
    'E49' = '  +2.48%  '
    'D50' = '0.07309'
    'E50' = '  +1.18%  '
    'D51' = '79.76'
    'E51' = '  +5.14%  '
}

foreach ($addr in $cells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $cells[$addr]
    $rng.Style = "Normal"
}
